$d = $word.ActiveDocument

# 1) Append a manual line break to the end of the last paragraph's existing run.
#    (No apostrophes in this search/replace text, so Find/Replace's smart-quote
#    autocorrect is not a concern here, and "^l" cleanly becomes a same-run <w:br/>.)
$null = $d.Content.Find.Execute("remain the same.", $true, $false, $false, $false, $false, $true, 1, $false, "remain the same.^l", 2)

# 2) Insert a brand-new paragraph after the (now last) paragraph; it inherits the
#    ListParagraph style + numPr + rPr from the paragraph it follows.
$lastP = $d.Paragraphs.Last
$lastP.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Last

# 3) Assign the new paragraph's text directly (not via Find/Replace) so the
#    straight apostrophes in "isn't"/"wouldn't" are not auto-corrected into
#    curly quotes. Embed Chr(11) (manual line break) where the diff wants <w:br/>.
$lb = [char]11
$newP.Range.Text = "The method that could be used to generate data that isn't sorted is quicksort because of its partition. Each partition is a logical separation of the data. Furthermore, the final position of each element within the partition is within the partition. Therefore, the locality of a set of a data is the size of the partition. Using this information, you could take any array, even if it isn't sorted, and form a non-sorted array that has the locality condition." + $lb + $lb + "To create the array with a specified locality, call quicksort as normal until a partition size is less than or equal to the locality parameter. Quicksort wouldn't be called on this partition (of size <= d). This would act as a base case in the recursion. Because the locality of an element in the quicksort algorithm is the size of its partition, stopping when partitions are less than d would create data that possesses the locality condition."

Write-Output "done"
